$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.282.95'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '2.615.63'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.598'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.21%  '
$ws.Range("D14").Value = '3.015.89'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").Value = '2.618.04'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.920'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").Value = '46.538.86'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("E20").Value = '  +2.87%  '
$ws.Range("E21").Value = '  +3.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '293.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.01%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '39.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("E32").Value = '  -1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("E37").Value = '  +3.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("E39").Value = '  +7.01%  '
$ws.Range("E40").Value = '  +2.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0332'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.28%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.39%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").Value = '2.116.63'
$ws.Range("E46").Value = '  +3.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.203'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.00%  '
